$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 580.6286
$ws.Range("I15").Value = 580.6286
$ws.Range("K15").Value = 1741.8858
$ws.Range("M15").Value = -1572.8858
$ws.Range("H31").Value = 500
$ws.Range("I31").Value = 500
$ws.Range("K31").Value = 1500
$ws.Range("M31").Value = -1270
$ws.Range("H33").Value = 759.6
$ws.Range("I33").Value = 266
$ws.Range("J33").Value = 1500
$ws.Range("K33").Value = 266
$ws.Range("L33").Value = 1500
$ws.Range("M33").Value = -37
$ws.Range("N33").Value = -1958
$ws.Range("H40").Value = 6858.125
$ws.Range("I40").Value = 3955.3333
$ws.Range("K40").Value = 3955.3333
$ws.Range("M40").Value = -3780.3333
$ws.Range("H86").Value = 4924.6665
$ws.Range("I86").Value = 4887.5
$ws.Range("K86").Value = 4887.5
$ws.Range("M86").Value = -3764.5
$ws.Range("H87").Value = 74354.25
$ws.Range("J87").Value = 93806.336
$ws.Range("L87").Value = 93806.336
$ws.Range("N87").Value = -96302.336
$ws.Range("H89").Value = 4924.6665
$ws.Range("I89").Value = 4887.5
$ws.Range("K89").Value = 24437.5
$ws.Range("M89").Value = -18821.5
$ws.Range("H90").Value = 74354.25
$ws.Range("J90").Value = 93806.336
$ws.Range("L90").Value = 281419.008
$ws.Range("N90").Value = -293899.008
$ws.Range("H107").Value = 367.14285
$ws.Range("I107").Value = 192.9
$ws.Range("J107").Value = 802.75
$ws.Range("K107").Value = 192.9
$ws.Range("L107").Value = 802.75
$ws.Range("M107").Value = 1727.1
$ws.Range("N107").Value = -4642.75
$ws.Range("H137").Value = 3304.2144
$ws.Range("I137").Value = 2283
$ws.Range("J137").Value = 4665.8335
$ws.Range("K137").Value = 6849
$ws.Range("L137").Value = 13997.5005
$ws.Range("M137").Value = -4299
$ws.Range("N137").Value = -19097.5005
$ws.Range("H138").Value = 2408.4285
$ws.Range("I138").Value = 1900
$ws.Range("K138").Value = 5700
$ws.Range("M138").Value = -560

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H8").Value = 6833
$ws.Range("I8").Value = 7000
$ws.Range("J8").Value = 6749.5
$ws.Range("K8").Value = 7000
$ws.Range("L8").Value = 6749.5
$ws.Range("M8").Value = -6856
$ws.Range("N8").Value = -7037.5
$ws.Range("H32").Value = 10949.889
$ws.Range("I32").Value = 24518.334
$ws.Range("K32").Value = 24518.334
$ws.Range("M32").Value = -24231.334
$ws.Range("H35").Value = 2800
$ws.Range("I35").Value = 2800
$ws.Range("K35").Value = 2800
$ws.Range("M35").Value = -2394
$ws.Range("H61").Value = 5610.778
$ws.Range("I61").Value = 5285.2856
$ws.Range("J61").Value = 6750
$ws.Range("K61").Value = 5285.2856
$ws.Range("L61").Value = 6750
$ws.Range("M61").Value = -5073.2856
$ws.Range("N61").Value = -7174
$ws.Range("H122").Value = 1482.0344
$ws.Range("J122").Value = 2487.5
$ws.Range("L122").Value = 7462.5
$ws.Range("N122").Value = -12362.5
$ws.Range("H132").Value = 1631.3529
$ws.Range("I132").Value = 1631.3529
$ws.Range("K132").Value = 4894.0587
$ws.Range("M132").Value = -2364.0587
$ws.Range("H136").Value = 5610.778
$ws.Range("I136").Value = 5285.2856
$ws.Range("J136").Value = 6750
$ws.Range("K136").Value = 15855.8568
$ws.Range("L136").Value = 20250
$ws.Range("M136").Value = -13305.8568
$ws.Range("N136").Value = -25350

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1411.25
$ws.Range("I20").Value = 1183
$ws.Range("J20").Value = 3009
$ws.Range("K20").Value = 1183
$ws.Range("L20").Value = 3009
$ws.Range("M20").Value = -936
$ws.Range("N20").Value = -3503
$ws.Range("H37").Value = 1535.2
$ws.Range("I37").Value = 1494.125
$ws.Range("J37").Value = 1699.5
$ws.Range("K37").Value = 1494.125
$ws.Range("L37").Value = 1699.5
$ws.Range("M37").Value = -1357.125
$ws.Range("N37").Value = -1973.5
$ws.Range("H107").Value = 4914.0835
$ws.Range("I107").Value = 3643.3333
$ws.Range("J107").Value = 6184.8335
$ws.Range("K107").Value = 3643.3333
$ws.Range("L107").Value = 6184.8335
$ws.Range("M107").Value = -1723.3333
$ws.Range("N107").Value = -10024.8335
$ws.Range("H134").Value = 2698.6428
$ws.Range("J134").Value = 5000
$ws.Range("L134").Value = 15000
$ws.Range("N134").Value = -20070

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1300
$ws.Range("I16").Value = 1100
$ws.Range("J16").Value = 1500
$ws.Range("K16").Value = 1100
$ws.Range("L16").Value = 1500
$ws.Range("M16").Value = -813
$ws.Range("N16").Value = -2074
$ws.Range("H22").Value = 420.375
$ws.Range("I22").Value = 393.5
$ws.Range("K22").Value = 393.5
$ws.Range("M22").Value = -43.5
$ws.Range("H31").Value = 6938.778
$ws.Range("I31").Value = 3046.4
$ws.Range("J31").Value = 7823.409
$ws.Range("K31").Value = 3046.4
$ws.Range("L31").Value = 7823.409
$ws.Range("M31").Value = -2751.4
$ws.Range("N31").Value = -8413.409
$ws.Range("H34").Value = 6938.778
$ws.Range("I34").Value = 3046.4
$ws.Range("J34").Value = 7823.409
$ws.Range("K34").Value = 3046.4
$ws.Range("L34").Value = 7823.409
$ws.Range("M34").Value = -2844.4
$ws.Range("N34").Value = -8227.409
$ws.Range("H62").Value = 0
$ws.Range("I62").Value = 0
$ws.Range("J62").Value = 0
$ws.Range("K62").Value = 0
$ws.Range("L62").ClearContents()
$ws.Range("M62").ClearContents()
$ws.Range("N62").ClearContents()
$ws.Range("H65").Value = 0
$ws.Range("I65").Value = 0
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = 0
$ws.Range("L65").ClearContents()
$ws.Range("M65").ClearContents()
$ws.Range("N65").ClearContents()
$ws.Range("H113").Value = 1300
$ws.Range("I113").Value = 1100
$ws.Range("J113").Value = 1500
$ws.Range("K113").Value = 1100
$ws.Range("L113").Value = 1500
$ws.Range("M113").Value = 1070
$ws.Range("N113").Value = -5840
$ws.Range("H122").Value = 3611
$ws.Range("J122").Value = 2222
$ws.Range("L122").Value = 6666
$ws.Range("N122").Value = -11566

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 106897.945
$ws.Range("I4").Value = 502497.75
$ws.Range("J4").Value = 1404.6666
$ws.Range("K4").Value = 1507493.25
$ws.Range("L4").Value = 4213.9998
$ws.Range("M4").Value = -1507381.25
$ws.Range("N4").Value = -4437.9998
$ws.Range("H9").Value = 317.8
$ws.Range("I9").Value = 22
$ws.Range("K9").Value = 66
$ws.Range("M9").Value = 158
$ws.Range("H99").Value = 1816.6666
$ws.Range("I99").Value = 1225.5
$ws.Range("K99").Value = 3676.5
$ws.Range("M99").Value = -1430.5
$ws.Range("H113").Value = 519.7778
$ws.Range("I113").Value = 483.75
$ws.Range("K113").Value = 1451.25
$ws.Range("M113").Value = 718.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 8498
$ws.Range("I80").Value = 8498
$ws.Range("J80").Value = 0
$ws.Range("K80").Value = 8498
$ws.Range("L80").Value = 0
$ws.Range("M80").ClearContents()
$ws.Range("N80").ClearContents()
$ws.Range("H83").Value = 8498
$ws.Range("I83").Value = 8498
$ws.Range("J83").Value = 0
$ws.Range("K83").Value = 42490
$ws.Range("L83").Value = 0
$ws.Range("M83").ClearContents()
$ws.Range("N83").ClearContents()
$ws.Range("H132").Value = 2386.111
$ws.Range("I132").Value = 2145.6875
$ws.Range("J132").Value = 4309.5
$ws.Range("K132").Value = 6437.0625
$ws.Range("L132").Value = 12928.5
$ws.Range("M132").Value = -3907.0625
$ws.Range("N132").Value = -17988.5
$ws.Range("H135").Value = 95000
$ws.Range("J135").Value = 95000
$ws.Range("L135").Value = 95000
$ws.Range("N135").Value = -105140

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 10000.667
$ws.Range("I7").Value = 9998
$ws.Range("K7").Value = 9998
$ws.Range("M7").Value = -9886
$ws.Range("H61").Value = 3868.375
$ws.Range("I61").Value = 1848.75
$ws.Range("J61").Value = 5888
$ws.Range("K61").Value = 1848.75
$ws.Range("L61").Value = 5888
$ws.Range("M61").Value = -1646.75
$ws.Range("N61").Value = -6292
$ws.Range("H113").Value = 3868.375
$ws.Range("I113").Value = 1848.75
$ws.Range("J113").Value = 5888
$ws.Range("K113").Value = 1848.75
$ws.Range("L113").Value = 5888
$ws.Range("M113").Value = 321.25
$ws.Range("N113").Value = -10228
$ws.Range("H126").Value = 10000.667
$ws.Range("I126").Value = 9998
$ws.Range("K126").Value = 29994
$ws.Range("M126").Value = -27524
$ws.Range("H132").Value = 5198.6
$ws.Range("I132").Value = 4998.25
$ws.Range("K132").Value = 14994.75
$ws.Range("M132").Value = -12464.75

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H43").Value = 0
$ws.Range("I43").Value = 0
$ws.Range("K43").Value = 0
$ws.Range("M43").ClearContents()
$ws.Range("H100").Value = 792.53845
$ws.Range("I100").Value = 691.9167
$ws.Range("K100").Value = 1383.8334
$ws.Range("M100").Value = -842.8334
$ws.Range("H109").Value = 0
$ws.Range("J109").Value = 0
$ws.Range("L109").ClearContents()
$ws.Range("N109").ClearContents()
$ws.Range("H132").Value = 2038.5714
$ws.Range("I132").Value = 1710.5555
$ws.Range("J132").Value = 2629
$ws.Range("K132").Value = 5131.666499999999
$ws.Range("L132").Value = 7887
$ws.Range("M132").Value = -2601.666499999999
$ws.Range("N132").Value = -12947
